$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("295k")

# ---- Fill in the missing data for the existing "good fit" block (rows 38-42) ----
$ws.Range("B38").Value = 78
$ws.Range("C38").Value = 58
$ws.Range("D38").Value = 20
$ws.Range("E38").Value = 11
$ws.Range("F38").Value = 9

$ws.Range("C39").Value = 7
$ws.Range("D39").Value = 0
$ws.Range("E39").Value = 0

$ws.Range("C40").Value = 49
$ws.Range("D40").Value = 2
$ws.Range("E40").Value = 1
$ws.Range("F40").Value = 2

$ws.Range("C41").Value = 2
$ws.Range("D41").Value = 11
$ws.Range("E41").Value = 9
$ws.Range("F41").Value = 1

$ws.Range("C42").Value = 0
$ws.Range("D42").Value = 8
$ws.Range("E42").Value = 1
$ws.Range("F42").Value = 6

# ---- Add a new "good fit 3 signal" block (rows 45-50), mirroring the header layout ----
$ws.Range("B45").Value = "all"
$ws.Range("C45").Value = "f1_good"
$ws.Range("D45").Value = "f1_bad"
$ws.Range("E45").Value = "f2_good"
$ws.Range("F45").Value = "f2_bad"
$ws.Range("G45").Value = "f3_good"
$ws.Range("H45").Value = "f3_bad"
$ws.Range("J45").Value = "chi2_per_dof_th"

$ws.Range("A46").Value = "всего"
$ws.Range("B46").Value = 78
$ws.Range("C46").Value = 58
$ws.Range("D46").Value = 20
$ws.Range("E46").Value = 10
$ws.Range("F46").Value = 10
$ws.Range("G46").Value = 7
$ws.Range("H46").Value = 3
$ws.Range("J46").Value = 5

$ws.Range("A47").Value = "шумы"
$ws.Range("A48").Value = "одиночные"
$ws.Range("A49").Value = "двойные"
$ws.Range("A50").Value = "тройные"

# ---- Restore the view state (active cell / top-left scroll position) ----
$ws.Range("I31").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
